$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Update existing rows 176-185 (Vega Monumental Concepción - Cereza)
# with the new weekly data per the commit.
# ---------------------------------------------------------------------

# Row 176
$ws.Range("D176").Value = 44931
$ws.Range("K176").Value = 'Lapins'
$ws.Range("L176").Value = 'Primera'
$ws.Range("M176").Value = 100
$ws.Range("N176").Value = 5000
$ws.Range("O176").Value = 6000
$ws.Range("P176").Value = 5500
$ws.Range("Q176").Value = '$/bandeja 10 kilos'
$ws.Range("R176").Value = "Región de O'Higgins"
$ws.Range("S176").Value = 550

# Row 177
$ws.Range("D177").Value = 44931
$ws.Range("K177").Value = 'Lapins'
$ws.Range("L177").Value = 'Segunda'
$ws.Range("M177").Value = 50
$ws.Range("N177").Value = 4000
$ws.Range("O177").Value = 4000
$ws.Range("P177").Value = 4000
$ws.Range("Q177").Value = '$/bandeja 10 kilos'
$ws.Range("R177").Value = "Región de O'Higgins"
$ws.Range("S177").Value = 400

# Row 178
$ws.Range("D178").Value = 44890
$ws.Range("K178").Value = 'Santina'
$ws.Range("L178").Value = 'Primera'
$ws.Range("M178").Value = 200
$ws.Range("N178").Value = 10000
$ws.Range("O178").Value = 11000
$ws.Range("P178").Value = 10500
$ws.Range("Q178").Value = '$/caja 10 kilos'
$ws.Range("R178").Value = 'Región de Ñuble'
$ws.Range("S178").Value = 1050

# Row 179
$ws.Range("D179").Value = 44890
$ws.Range("K179").Value = 'Santina'
$ws.Range("L179").Value = 'Segunda'
$ws.Range("M179").Value = 100
$ws.Range("N179").Value = 8000
$ws.Range("O179").Value = 8000
$ws.Range("P179").Value = 8000
$ws.Range("Q179").Value = '$/caja 10 kilos'
$ws.Range("R179").Value = 'Región de Ñuble'
$ws.Range("S179").Value = 800

# Row 180
$ws.Range("D180").Value = 44580
$ws.Range("K180").Value = 'Lapins'
$ws.Range("L180").Value = 'Primera'
$ws.Range("M180").Value = 220
$ws.Range("N180").Value = 4500
$ws.Range("O180").Value = 5000
$ws.Range("P180").Value = 4773
$ws.Range("Q180").Value = '$/bandeja 10 kilos'
$ws.Range("R180").Value = 'Provincia de Curicó'
$ws.Range("S180").Value = 477

# Row 181
$ws.Range("D181").Value = 44565
$ws.Range("K181").Value = 'Lapins'
$ws.Range("L181").Value = 'Primera'
$ws.Range("M181").Value = 250
$ws.Range("N181").Value = 4000
$ws.Range("O181").Value = 5000
$ws.Range("P181").Value = 4400
$ws.Range("Q181").Value = '$/caja 10 kilos'
$ws.Range("R181").Value = 'Provincia de Curicó'
$ws.Range("S181").Value = 440

# Row 182
$ws.Range("D182").Value = 44565
$ws.Range("K182").Value = 'Rainier'
$ws.Range("L182").Value = 'Primera'
$ws.Range("M182").Value = 220
$ws.Range("N182").Value = 5000
$ws.Range("O182").Value = 5500
$ws.Range("P182").Value = 5227
$ws.Range("Q182").Value = '$/bandeja 10 kilos'
$ws.Range("R182").Value = 'Provincia de Curicó'
$ws.Range("S182").Value = 523

# Row 183
$ws.Range("D183").Value = 44589
$ws.Range("K183").Value = 'Santina'
$ws.Range("L183").Value = 'Primera'
$ws.Range("M183").Value = 160
$ws.Range("N183").Value = 7500
$ws.Range("O183").Value = 8000
$ws.Range("P183").Value = 7750
$ws.Range("Q183").Value = '$/bandeja 10 kilos'
$ws.Range("R183").Value = 'Provincia de Curicó'
$ws.Range("S183").Value = 775

# Row 184 (only variety changes)
$ws.Range("K184").Value = 'Lapins'

# Row 185 (only variety changes)
$ws.Range("K185").Value = 'Lapins'

# ---------------------------------------------------------------------
# Append two brand-new rows (186, 187) at the end of the table, carrying
# forward the data that used to live in rows 184/185 before the variety
# there was changed to Lapins (this becomes a fresh Santina record).
# ---------------------------------------------------------------------

# Row 186
$ws.Range("A186").Value = 11
$ws.Range("B186").Value = 'Vega Monumental Concepción'
$ws.Range("C186").Value = 'Bíobío'
$ws.Range("D186").Value = 44911
$ws.Range("D186").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E186").Value = 8
$ws.Range("F186").Value = 'Fruta'
$ws.Range("G186").Value = 100103
$ws.Range("H186").Value = 'Frutos de hueso (carozo)'
$ws.Range("I186").Value = 100103001
$ws.Range("J186").Value = 'Cereza'
$ws.Range("K186").Value = 'Santina'
$ws.Range("L186").Value = 'Primera'
$ws.Range("M186").Value = 100
$ws.Range("N186").Value = 5000
$ws.Range("O186").Value = 6000
$ws.Range("P186").Value = 5500
$ws.Range("Q186").Value = '$/caja 10 kilos'
$ws.Range("R186").Value = 'Región de Ñuble'
$ws.Range("S186").Value = 550
$ws.Range("T186").Value = 10

# Row 187
$ws.Range("A187").Value = 11
$ws.Range("B187").Value = 'Vega Monumental Concepción'
$ws.Range("C187").Value = 'Bíobío'
$ws.Range("D187").Value = 44911
$ws.Range("D187").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E187").Value = 8
$ws.Range("F187").Value = 'Fruta'
$ws.Range("G187").Value = 100103
$ws.Range("H187").Value = 'Frutos de hueso (carozo)'
$ws.Range("I187").Value = 100103001
$ws.Range("J187").Value = 'Cereza'
$ws.Range("K187").Value = 'Santina'
$ws.Range("L187").Value = 'Segunda'
$ws.Range("M187").Value = 50
$ws.Range("N187").Value = 4000
$ws.Range("O187").Value = 4000
$ws.Range("P187").Value = 4000
$ws.Range("Q187").Value = '$/caja 10 kilos'
$ws.Range("R187").Value = 'Región de Ñuble'
$ws.Range("S187").Value = 400
$ws.Range("T187").Value = 10
